$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right mark 5 -> 4, Wrong mark -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total 80 -> 64, Wrong total -5 -> -10, summary 80/140 -> 54/112
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "54 / 112"
